$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 1071 (existing rows 1071-1148 shift down to 1074-1151).
$ws.Range("A1071:A1073").EntireRow.Insert()

# Common values shared by all rows in this sheet (market / region / product taxonomy).
$mercadoId = 10
$mercado = "Vega Modelo de Temuco"
$region = "La Araucanía"
$codreg = 9
$tipo = "Fruta"
$productoId = 100104
$producto = "Frutos de pepita"
$categoriaId = 100104005
$categoria = "Pera"
$origen = "Región de O'Higgins"

$fecha = 45013
$unidad = "$/bandeja 18 kilos granel"
$calidad = "Primera"
$kgUnidad = 18

# Row 1071: Forelle
$r = 1071
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Forelle"
$ws.Cells.Item($r, 12).Value = $calidad
$ws.Cells.Item($r, 13).Value = 110
$ws.Cells.Item($r, 14).Value = 12000
$ws.Cells.Item($r, 15).Value = 12000
$ws.Cells.Item($r, 16).Value = 12000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 667
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 1072: Packham's Triumph
$r = 1072
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Packham's Triumph"
$ws.Cells.Item($r, 12).Value = $calidad
$ws.Cells.Item($r, 13).Value = 120
$ws.Cells.Item($r, 14).Value = 12000
$ws.Cells.Item($r, 15).Value = 12000
$ws.Cells.Item($r, 16).Value = 12000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 667
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 1073: Winter Nelis
$r = 1073
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Winter Nelis"
$ws.Cells.Item($r, 12).Value = $calidad
$ws.Cells.Item($r, 13).Value = 95
$ws.Cells.Item($r, 14).Value = 12000
$ws.Cells.Item($r, 15).Value = 12000
$ws.Cells.Item($r, 16).Value = 12000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 667
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Ensure date formatting is preserved for column D on the new rows.
$ws.Range("D1071:D1073").NumberFormat = "YYYY-MM-DD HH:MM:SS"
